# RPA datasets push 2024-08-08
#
# Update confirmed IPO offering price (확정공모가, column D) and the
# fundraising amount (공모금액(백만), column E) for three listings whose
# bookbuilding has since concluded:
#   row 11 - 케이쓰리아이         : D "-"  -> "15500", E "17500" -> "22351"
#   row 12 - 전진건설로봇(유가)   : D "-"  -> "16500", E "42471" -> "50781"
#   row 14 - 넥스트바이오메디컬   : D "-"  -> "29000", E "24000" -> "29000"
#
# These columns store every value as text (shared-string) cells, so the
# new numeric-looking values must also be written as text rather than
# numbers. A leading apostrophe forces Excel to keep the entry as text,
# and resetting the cell style back to "Normal" afterwards removes the
# quote-prefix formatting flag so the cell is left exactly as it would be
# for any other plain text cell in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "'15500"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "'22351"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'16500"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "'50781"
$ws.Range("E12").Style = "Normal"

$ws.Range("D14").Value = "'29000"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "'29000"
$ws.Range("E14").Style = "Normal"
